$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("small_anchor")
for ($w = 5; $w -le 20; $w++) {
    $ws.Columns.Item($w - 4).ColumnWidth = $w
}
Write-Output "done"
